# Add data for 2022-01-05: update workbook/sheet title and the December /
# Total rows of the carjacking-by-month-yoy table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and its tab) from "Through 2021-12-27" to "Through 2021-12-28"
$ws.Name = "Through 2021-12-28"

# Update the label for the December row to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-28)"

# Update the December row (row 13) values for each year column (B:H)
$ws.Range("B13").Value = 42
$ws.Range("C13").Value = 91
$ws.Range("D13").Value = 108
$ws.Range("E13").Value = 68
$ws.Range("F13").Value = 57
$ws.Range("G13").Value = 133
$ws.Range("H13").Value = 170

# Update the Total row (row 14) values for each year column (B:H)
$ws.Range("B14").Value = 333
$ws.Range("C14").Value = 654
$ws.Range("D14").Value = 929
$ws.Range("E14").Value = 750
$ws.Range("F14").Value = 591
$ws.Range("G14").Value = 1397
$ws.Range("H14").Value = 1813
